$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FBS")
$ws2 = $wb.Worksheets.Item("Other")

# --- Update Timestamp shared string (column AK on FBS, rows 2-54) ---
$ws1.Range("AK2:AK54").Value = "2024-10-20T16:21:43.793690"

# --- FBS sheet per-row updates (betting odds / move columns) ---
# Row 2
$ws1.Range("Y2").Value = 37.5
$ws1.Range("Z2").Value = -110
$ws1.Range("AE2").Value = 0.0273972602739726

# Row 3
$ws1.Range("Y3").Value = 66.5
$ws1.Range("Z3").Value = -110
$ws1.Range("AE3").Value = 0.0310077519379845

# Row 4
$ws1.Range("W4").Value = 64.5
$ws1.Range("X4").Value = -105
$ws1.Range("Y4").Value = 64.5
$ws1.Range("Z4").Value = -105
$ws1.Range("AE4").Value = 0

# Row 5
$ws1.Range("Y5").Value = 45.5
$ws1.Range("Z5").Value = -115
$ws1.Range("AE5").Value = -0.02150537634408602

# Row 6
$ws1.Range("W6").Value = 49.5
$ws1.Range("X6").Value = -110
$ws1.Range("Y6").Value = 49.5
$ws1.Range("Z6").Value = -110
$ws1.Range("AE6").Value = 0

# Row 7
$ws1.Range("W7").Value = 59.5
$ws1.Range("X7").Value = -110
$ws1.Range("Y7").Value = 59.5
$ws1.Range("Z7").Value = -110
$ws1.Range("AE7").Value = 0

# Row 9
$ws1.Range("Y9").Value = 64.5
$ws1.Range("Z9").Value = -115
$ws1.Range("AE9").Value = 0.08403361344537816

# Row 10
$ws1.Range("Y10").Value = 51.5
$ws1.Range("Z10").Value = -105
$ws1.Range("AE10").Value = -0.01904761904761905

# Row 11
$ws1.Range("W11").Value = 46.5
$ws1.Range("X11").Value = -110
$ws1.Range("Y11").Value = 46.5
$ws1.Range("Z11").Value = -110
$ws1.Range("AE11").Value = 0
$ws1.Range("AA11").ClearContents()
$ws1.Range("AB11").ClearContents()

# Row 12
$ws1.Range("W12").Value = 65.5
$ws1.Range("X12").Value = -110
$ws1.Range("Y12").Value = 65.5
$ws1.Range("Z12").Value = -110
$ws1.Range("AE12").Value = 0

# Row 14
$ws1.Range("Y14").Value = 55.5
$ws1.Range("AE14").Value = 0

# Row 15
$ws1.Range("Q15").Value = "N"
$ws1.Range("Y15").Value = 56.5
$ws1.Range("AB15").Value = -14
$ws1.Range("AE15").Value = 0
$ws1.Range("AF15").Value = -0.5

# Row 16
$ws1.Range("Y16").Value = 58.5
$ws1.Range("AE16").Value = -0.01680672268907563

# Row 19
$ws1.Range("W19").Value = 57.5
$ws1.Range("X19").Value = -110
$ws1.Range("Y19").Value = 57.5
$ws1.Range("Z19").Value = -110
$ws1.Range("AE19").Value = 0

# Row 22
$ws1.Range("Y22").Value = 63.5
$ws1.Range("Z22").Value = -105
$ws1.Range("AE22").Value = 0

# Row 23
$ws1.Range("W23").Value = 54.5
$ws1.Range("X23").Value = -105
$ws1.Range("Y23").Value = 54.5
$ws1.Range("Z23").Value = -105
$ws1.Range("AE23").Value = 0

# Row 24
$ws1.Range("Y24").Value = 60.5
$ws1.Range("AE24").Value = 0.01680672268907563

# Row 25
$ws1.Range("Y25").Value = 53.5
$ws1.Range("Z25").Value = -110
$ws1.Range("AE25").Value = 0

# Row 27
$ws1.Range("Y27").Value = 60.5
$ws1.Range("Z27").Value = -115
$ws1.Range("AE27").Value = 0.03418803418803419

# Row 29
$ws1.Range("W29").Value = 46.5
$ws1.Range("X29").Value = -110
$ws1.Range("Y29").Value = 46.5
$ws1.Range("Z29").Value = -110
$ws1.Range("AE29").Value = 0

# Row 32
$ws1.Range("W32").Value = 63.5
$ws1.Range("X32").Value = -115
$ws1.Range("Y32").Value = 63.5
$ws1.Range("Z32").Value = -115
$ws1.Range("AE32").Value = 0

# Row 33
$ws1.Range("W33").Value = 48.5
$ws1.Range("X33").Value = -115
$ws1.Range("Y33").Value = 48.5
$ws1.Range("Z33").Value = -115
$ws1.Range("AB33").Value = -7
$ws1.Range("AE33").Value = 0
$ws1.Range("AF33").Value = -0.5

# Row 34
$ws1.Range("W34").Value = 47.5
$ws1.Range("X34").Value = -110
$ws1.Range("Y34").Value = 47.5
$ws1.Range("Z34").Value = -110
$ws1.Range("AE34").Value = 0

# Row 35
$ws1.Range("Y35").Value = 54.5
$ws1.Range("AE35").Value = -0.01801801801801802

# Row 36
$ws1.Range("W36").Value = 51.5
$ws1.Range("X36").Value = -105
$ws1.Range("Y36").Value = 51.5
$ws1.Range("Z36").Value = -105
$ws1.Range("AE36").Value = 0

# Row 37
$ws1.Range("W37").Value = 54.5
$ws1.Range("X37").Value = -110
$ws1.Range("Y37").Value = 54.5
$ws1.Range("Z37").Value = -110
$ws1.Range("AE37").Value = 0

# Row 38
$ws1.Range("W38").Value = 46.5
$ws1.Range("X38").Value = -115
$ws1.Range("Y38").Value = 46.5
$ws1.Range("Z38").Value = -115
$ws1.Range("AE38").Value = 0

# Row 39
$ws1.Range("Z39").Value = -110

# Row 40
$ws1.Range("Q40").Value = "SSE"
$ws1.Range("W40").Value = 55.5
$ws1.Range("X40").Value = -110
$ws1.Range("Y40").Value = 55.5
$ws1.Range("Z40").Value = -110
$ws1.Range("AE40").Value = 0

# Row 42
$ws1.Range("W42").Value = 48.5
$ws1.Range("X42").Value = -110
$ws1.Range("Y42").Value = 48.5
$ws1.Range("Z42").Value = -110
$ws1.Range("AE42").Value = 0

# Row 43
$ws1.Range("Q43").Value = "N"
$ws1.Range("W43").Value = 53.5
$ws1.Range("X43").Value = -115
$ws1.Range("Y43").Value = 53.5
$ws1.Range("Z43").Value = -115
$ws1.Range("AE43").Value = 0

# Row 45
$ws1.Range("W45").Value = 56.5
$ws1.Range("X45").Value = -115
$ws1.Range("Y45").Value = 56.5
$ws1.Range("Z45").Value = -115
$ws1.Range("AE45").Value = 0

# Row 46
$ws1.Range("Y46").Value = 41.5
$ws1.Range("Z46").Value = -110
$ws1.Range("AE46").Value = 0.05063291139240506

# Row 50
$ws1.Range("Y50").Value = 48.5
$ws1.Range("AE50").Value = -0.0396039603960396

# Row 51
$ws1.Range("W51").Value = 54.5
$ws1.Range("X51").Value = -110
$ws1.Range("Y51").Value = 54.5
$ws1.Range("Z51").Value = -110
$ws1.Range("AE51").Value = 0

# Row 52
$ws1.Range("Y52").Value = 56.5
$ws1.Range("Z52").Value = -110
$ws1.Range("AE52").Value = -0.03418803418803419

# Row 53
$ws1.Range("W53").Value = 56.5
$ws1.Range("X53").Value = -115
$ws1.Range("Y53").Value = 56.5
$ws1.Range("Z53").Value = -115
$ws1.Range("AE53").Value = 0

# Row 54
$ws1.Range("W54").Value = 47.5
$ws1.Range("X54").Value = -115
$ws1.Range("Y54").Value = 47.5
$ws1.Range("Z54").Value = -115
$ws1.Range("AE54").Value = 0

# --- Other sheet wind_dir_fg (S column) updates ---
$ws2.Range("S16").Value = "E"
$ws2.Range("S36").Value = "SSE"
